$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last data row (old row 6); the table shrinks from 5 data rows to 4.
$ws.Rows(6).Delete()

# Row 2: product code/description updated to what used to be row 3's item;
# quantity (D) stays the same, rest (E/F) unchanged.
$ws.Range("A2").Value = "20135336"
$ws.Range("B2").Value = "CMORY BITE STR LC120"

# Row 3: now holds the item that used to be on row 4 (D/E/F), but keeps the
# product code/description from row 2's new item.
$ws.Range("A3").Value = "20135336"
$ws.Range("B3").Value = "CMORY BITE STR LC120"
$ws.Range("D3").Value = "7"
$ws.Range("E3").Value = "16"
$ws.Range("F3").Value = "RT,(E-1B)"

# Row 4: now holds what used to be on row 5.
$ws.Range("A4").Value = "20087415"
$ws.Range("B4").Value = "SG CHKN.NUG ALPH 200"
$ws.Range("D4").Value = "7"
$ws.Range("E4").Value = "22"
$ws.Range("F4").Value = "RT,(E-1B)"

# Row 5: now holds what used to be on row 6 (before it was deleted above).
$ws.Range("A5").Value = "20093522"
$ws.Range("B5").Value = "SG SPCY CHKN STRP250"
$ws.Range("D5").Value = "7"
$ws.Range("E5").Value = "24"
$ws.Range("F5").Value = "RT,(E-1B)"
